# Add a new "right_border" style row (row 37) to Sheet1, mirroring the
# existing pattern of named-style / value pairs in column A / B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label in column A with a thin right border applied, and matching
# value in column B (consistent with the other rows in the sheet).
$ws.Range("A37").Value = "right_border"
$ws.Range("B37").Value = 11

# Apply a thin border on the right edge of A37 only (new borderId/cellXf).
# xlEdgeRight = 10, xlThin = 2
$ws.Range("A37").Borders.Item(10).LineStyle = 1
$ws.Range("A37").Borders.Item(10).Weight = 2

# Mirror the view state changes captured in the diff: select B38 (the cell
# just below the new data).
$ws.Range("B38").Select()
